$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.619.55"
$ws.Range("E2").Value = "  +2.18%  "
$ws.Range("D3").Value = "1.892.23"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'244.15"
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").Value = "'0.4961"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.2957"
$ws.Range("E8").Value = "  +1.81%  "
$ws.Range("D9").Value = "'0.06817"
$ws.Range("E9").Value = "  +3.70%  "
$ws.Range("D10").Value = "1.892.37"
$ws.Range("E10").Value = "  +1.04%  "
$ws.Range("D11").Value = "'17.08"
$ws.Range("E11").Value = "  +2.34%  "
$ws.Range("D12").Value = "'0.07334"
$ws.Range("E12").Value = "  +2.23%  "
$ws.Range("D13").Value = "'91.23"
$ws.Range("E13").Value = "  +6.10%  "
$ws.Range("D14").Value = "'5.096"
$ws.Range("E14").Value = "  +5.29%  "
$ws.Range("D15").Value = "'0.6742"
$ws.Range("E15").Value = "  +2.41%  "
$ws.Range("D16").Value = "30.626.00"
$ws.Range("E16").Value = "  +2.28%  "
$ws.Range("D17").Value = "'0.000007922"
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").Value = "'13.26"
$ws.Range("E19").Value = "  +4.77%  "
$ws.Range("D20").Value = "2.136.91"
$ws.Range("E20").Value = "  +1.08%  "
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("D22").Value = "'4.866"
$ws.Range("E22").Value = "  +3.00%  "
$ws.Range("D23").Value = "'178.33"
$ws.Range("E23").Value = "  +32.66%  "
$ws.Range("D24").Value = "'6.063"
$ws.Range("E24").Value = "  +8.71%  "
$ws.Range("D25").Value = "'9.288"
$ws.Range("E25").Value = "  +2.83%  "
$ws.Range("D26").Value = "'154.23"
$ws.Range("E26").Value = "  +2.88%  "
$ws.Range("D27").Value = "'18.78"
$ws.Range("E27").Value = "  +12.72%  "
$ws.Range("E28").Value = "  +1.66%  "
$ws.Range("D29").Value = "'1.387"
$ws.Range("E29").Value = "  +1.26%  "
$ws.Range("D30").Value = "'4.333"
$ws.Range("E30").Value = "  +4.54%  "
$ws.Range("D31").Value = "'0.08937"
$ws.Range("E31").Value = "  +2.96%  "
$ws.Range("D32").Value = "'4.038"
$ws.Range("E32").Value = "  +2.92%  "
$ws.Range("D33").Value = "'0.05203"
$ws.Range("E33").Value = "  +3.82%  "
$ws.Range("D34").Value = "'0.7395"
$ws.Range("E34").Value = "  +5.86%  "
$ws.Range("D35").Value = "'1.133"
$ws.Range("E35").Value = "  +3.90%  "
$ws.Range("E36").Value = "  +0.70%  "
$ws.Range("D37").Value = "'0.01873"
$ws.Range("E37").Value = "  +10.56%  "
$ws.Range("D38").Value = "'2.704"
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("D39").Value = "'2.168"
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("D40").Value = "'0.9333"
$ws.Range("E40").Value = "  +0.80%  "
$ws.Range("D41").Value = "'0.4359"
$ws.Range("E41").Value = "  +4.53%  "
$ws.Range("D42").Value = "'106.05"
$ws.Range("E42").Value = "  +4.51%  "
$ws.Range("E44").Value = "  +0.32%  "
$ws.Range("D45").Value = "'7.645"
$ws.Range("E45").Value = "  +3.75%  "
$ws.Range("D46").Value = "'0.1353"
$ws.Range("E46").Value = "  +7.92%  "
$ws.Range("D47").Value = "'0.05846"
$ws.Range("E47").Value = "  +3.66%  "
$ws.Range("D48").Value = "'33.42"
$ws.Range("E48").Value = "  +3.20%  "
$ws.Range("D49").Value = "'0.3892"
$ws.Range("E49").Value = "  +5.66%  "
$ws.Range("D50").Value = "'8.533"
$ws.Range("E50").Value = "  +5.66%  "
$ws.Range("D51").Value = "'1.378"
$ws.Range("E51").Value = "  +3.75%  "
